$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.250.26'
$ws.Range('D2').NumberFormat = 'General'
$ws.Range('E2').Value = '  +0.38%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.863.60'
$ws.Range('D3').NumberFormat = 'General'
$ws.Range('E3').Value = '  +0.12%  '

# Row 4
$ws.Range('E4').Value = '  +0.10%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '237.01'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  +1.33%  '

# Row 6
$ws.Range('E6').Value = '  +0.07%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4681'
$ws.Range('D7').NumberFormat = 'General'
$ws.Range('E7').Value = '  +0.50%  '

# Row 8
$ws.Range('E8').Value = '  +1.52%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06559'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('E9').Value = '  +0.16%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '22.18'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('E10').Value = '  +10.44%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07894'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('E11').Value = '  +1.04%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '98.04'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('E12').Value = '  +1.94%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.873.32'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('E13').Value = '  +0.68%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.193'
$ws.Range('D14').NumberFormat = 'General'

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6828'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('E15').Value = '  +2.04%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '277.60'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('E16').Value = '  -0.92%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.249.25'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('E17').Value = '  +0.31%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.63'
$ws.Range('D18').NumberFormat = 'General'

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.001'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').Value = '  +0.08%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007342'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').Value = '  +1.39%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.114.27'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').Value = '  +0.82%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.361'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  -2.04%  '

# Row 23
$ws.Range('E23').Value = '  +0.16%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.194'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  +0.95%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '168.11'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E25').Value = '  +1.92%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.244'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').Value = '  -0.74%  '

# Row 27
$ws.Range('E27').Value = '  +1.12%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.956'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('E28').Value = '  +2.91%  '

# Row 29
$ws.Range('E29').Value = '  +3.14%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.09848'
$ws.Range('D30').NumberFormat = 'General'

# Row 31
$ws.Range('E31').Value = '  -0.20%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.483'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('E32').Value = '  +1.01%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.074'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('E33').Value = '  -0.54%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04747'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('E34').Value = '  +2.07%  '

# Row 35
$ws.Range('E35').Value = '  +4.07%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7047'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  +0.69%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.710'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').Value = '  +0.03%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01879'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').Value = '  +1.48%  '

# Row 39
$ws.Range('E39').Value = '  +4.09%  '

# Row 40
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.290'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').Value = '  -0.04%  '

# Row 41
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '75.59'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('E41').Value = '  +4.01%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.957'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('E42').Value = '  +1.96%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8554'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').Value = '  +0.39%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.4179'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('E44').Value = '  +0.52%  '

# Row 45
$ws.Range('E45').Value = '  +0.04%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '103.50'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').Value = '  +0.09%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.223'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('E47').Value = '  +0.80%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '946.90'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  -4.35%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.244'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('E49').Value = '  +0.84%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '34.26'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('E50').Value = '  +0.40%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05643'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('E51').Value = '  +0.10%  '
